# osmanbo add stage-3 B1-B2
# Fill in the TxHash evidence rows for sheets "B1" and "B2", matching the
# pattern already used for the "A1".."A20" sheets, then restore the
# selection / active-sheet state to what it was left in afterwards.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# B1: first Interchain NFT-Transfer TxHash / Internal Transfer TxHash
# ---------------------------------------------------------------------
$wsB1 = $wb.Worksheets.Item("B1")
$wsB1.Range("A2").Value = "B59D759CFC2B3C2467393581B0279F6D59B7A5E072B382822337C16AC1EF6355"
$wsB1.Range("A3").Value = "A9FBE0099069CFA335599822C53EA91212E0E282AC286CB91FB58C25E45C3C53"

# Widen column A so the long hash is readable, and wrap the text in A3.
$wsB1.Columns.Item(1).ColumnWidth = 99.8
$wsB1.Range("A3").WrapText = $true

# ---------------------------------------------------------------------
# B2: first Interchain NFT-Transfer TxHash / Internal Transfer TxHash
# ---------------------------------------------------------------------
$wsB2 = $wb.Worksheets.Item("B2")
$wsB2.Range("A2").Value = "9A8E8DD23F9B39130FE14E70BC9A85398D76DFDAAC2E5D5969B4188D305E3AC2"
$wsB2.Range("A3").Value = "34F863AD81C8009FFFB0FD67EB2464852BB887EB3C64862A56A0F37099FAB5E4"

$wsB2.Columns.Item(1).ColumnWidth = 84.8

# ---------------------------------------------------------------------
# Restore / update the view state: Info, A20, B1 each keep their own
# selection; B2 ends up the active sheet/tab (tabSelected + activeTab).
# ---------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Activate()
$wsInfo.Range("D2").Select()

$wsA20 = $wb.Worksheets.Item("A20")
$wsA20.Activate()
$wsA20.Range("E22").Select()

$wsB1.Activate()
$wsB1.Range("A3").Select()

$wsB2.Activate()
$wsB2.Range("C10").Select()
